# "Added more C fibre data." -- update the Exp 19 baseline source file name
# and refresh the ramp-position figures for rows 2-3 (Exp 27 unit 1 / unit 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full file name for the Exp 19 baseline row (A2) now points at the new .mat file.
$ws.Range("A2").Value = "/home/daniel/Spike Data/Matlab files/Exp 19 baseline data new.mat"

# Row 2 (Exp 27 unit 1): updated channel + ramp position data.
$ws.Range("B2").Value = 2
$ws.Range("D2").Value = 25470
$ws.Range("E2").Value = 27430
$ws.Range("F2").Value = 119400
$ws.Range("G2").Value = 121200
$ws.Range("H2").Value = 199600
$ws.Range("I2").Value = 201400

# Row 3 (Exp 27 unit 2): updated start/end ramp figures for the second ramp.
$ws.Range("F3").Value = 2976
$ws.Range("G3").Value = 5271

# Window chrome: sheet-tab split ratio and the last active selection.
$excel.ActiveWindow.TabRatio = 24
$ws.Range("F11").Select()
